$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "gender"

$genders = @("m","f","m","m","f","m","f","m","f","m","m","m","m","m","m","f","f","f","f","f","m","m","m","f","m","m","f","f","f","m","m","f","f","f","m","m")

for ($i = 0; $i -lt $genders.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $genders[$i]
}

$ws.Columns.Item(7).ColumnWidth = 16.5

$ws.Range("J9").Select()
